{"js": "// Update the worksheet date and all the multiplication problems to the\n// next day's values, per the commit's regenerated content.\nconst replacements = [\n  [\"2024-09-10 Tuesday\", \"2024-09-11 Wednesday\"],\n  [\"430\u00d76=\", \"648\u00d75=\"],\n  [\"891\u00d79=\", \"780\u00d72=\"],\n  [\"813\u00d73=\", \"279\u00d73=\"],\n  [\"961\u00d77=\", \"661\u00d75=\"],\n  [\"878\u00d78=\", \"815\u00d76=\"],\n  [\"302\u00d72=\", \"382\u00d76=\"],\n  [\"109\u00d79=\", \"142\u00d78=\"],\n  [\"783\u00d73=\", \"390\u00d75=\"],\n  [\"689\u00d79=\", \"326\u00d73=\"],\n  [\"582\u00d74=\", \"545\u00d79=\"],\n  [\"905\u00d78=\", \"735\u00d73=\"],\n  [\"174\u00d72=\", \"816\u00d76=\"],\n  [\"293\u00d77=\", \"277\u00d78=\"],\n  [\"970\u00d78=\", \"972\u00d74=\"],\n  [\"408\u00d73=\", \"296\u00d79=\"],\n  [\"504\u00d79=\", \"354\u00d72=\"],\n  [\"657\u00d74=\", \"764\u00d72=\"],\n  [\"992\u00d79=\", \"120\u00d76=\"],\n  [\"349\u00d73=\", \"171\u00d79=\"],\n  [\"798\u00d77=\", \"607\u00d73=\"],\n  [\"430\u00d73=\", \"962\u00d78=\"],\n  [\"742\u00d77=\", \"842\u00d73=\"],\n  [\"422\u00d75=\", \"736\u00d79=\"],\n  [\"153\u00d73=\", \"772\u00d75=\"],\n  [\"157\u00d73=\", \"624\u00d79=\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all the multiplication problems to the\n# next day's values, per the commit's regenerated content.\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n  \"2024-09-10 Tuesday\" = \"2024-09-11 Wednesday\"\n  \"430\u00d76=\" = \"648\u00d75=\"\n  \"891\u00d79=\" = \"780\u00d72=\"\n  \"813\u00d73=\" = \"279\u00d73=\"\n  \"961\u00d77=\" = \"661\u00d75=\"\n  \"878\u00d78=\" = \"815\u00d76=\"\n  \"302\u00d72=\" = \"382\u00d76=\"\n  \"109\u00d79=\" = \"142\u00d78=\"\n  \"783\u00d73=\" = \"390\u00d75=\"\n  \"689\u00d79=\" = \"326\u00d73=\"\n  \"582\u00d74=\" = \"545\u00d79=\"\n  \"905\u00d78=\" = \"735\u00d73=\"\n  \"174\u00d72=\" = \"816\u00d76=\"\n  \"293\u00d77=\" = \"277\u00d78=\"\n  \"970\u00d78=\" = \"972\u00d74=\"\n  \"408\u00d73=\" = \"296\u00d79=\"\n  \"504\u00d79=\" = \"354\u00d72=\"\n  \"657\u00d74=\" = \"764\u00d72=\"\n  \"992\u00d79=\" = \"120\u00d76=\"\n  \"349\u00d73=\" = \"171\u00d79=\"\n  \"798\u00d77=\" = \"607\u00d73=\"\n  \"430\u00d73=\" = \"962\u00d78=\"\n  \"742\u00d77=\" = \"842\u00d73=\"\n  \"422\u00d75=\" = \"736\u00d79=\"\n  \"153\u00d73=\" = \"772\u00d75=\"\n  \"157\u00d73=\" = \"624\u00d79=\"\n}\n\nforeach ($find in $replacements.Keys) {\n  $replace = $replacements[$find]\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, \"wdReplaceAll\") | Out-Null\n}\n"}
